# Add a new event of blockway (row insert into SceneQuest Sheet1 table)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Insert a new blank row at position 64 (pushes existing rows 64-80 down
#    to 65-81, and inherits formatting/styles from the row above - row 63 -
#    which already matches what the new row needs: s=5 on I,K:P and s=24 on J).
$ws.Rows.Item(64).Insert()

# 2) Fill in the new row 64 with the "blockway" event data.
$ws.Range("A64").Value = 42010054
$ws.Range("B64").Value = "设卡"
$ws.Range("C64").Value = 1
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 2
$ws.Range("F64").Value = "blockway"
$ws.Range("G64").Value = "blockway"
$ws.Range("H64").Value = "blockway"
$ws.Range("J64").Value = "冒险"
$ws.Range("Q64").Value = "soldier"
$ws.Range("Y64").Value = 300
$ws.Range("AD64").Value = "ziyuandai(renlei)"
$ws.Range("AE64").Value = "dlsucaidai"
$ws.Range("AI64").Value = 150
$ws.Range("AJ64").Value = 100
$ws.Range("AK64").Value = 150
$ws.Range("AN64").Value = 40
$ws.Range("AO64").Value = 120
$ws.Range("AP64").Value = 50
$ws.Range("AQ64").Value = 50
$ws.Range("AR64").Value = 25

# 3) Re-assert G63 (string "starve") -- value unchanged but kept explicit
#    for clarity/consistency with the source edit.
$ws.Range("G63").Value = "starve"

# 4) Fix up the sheet view: drop the stale topLeftCell on the view, move the
#    frozen pane's scroll anchor down, and move the active selection.
$view = $ws.Application.ActiveWindow
$view.SplitRow = 3
$ws.Range("A46").Select()
$ws.Range("Z64").Select()

# 5) Extend the conditionalFormatting range that covers the data rows below
#    the header so it still reaches the bottom of the table (now row 81).
$ws.Range("B64:AW81").FormatConditions.Delete()
$cf = $ws.Range("B64:AW81").FormatConditions.Add(8, 0, "=LEN(TRIM(B64))=0")

# 6) Resize the table/list object so it (and its AutoFilter) cover the new
#    last row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:AW81"))
